$wb = $excel.ActiveWorkbook

# Germany sheet: selection changes to a "select all" style range (A1:XFD1048576)
$wsGermany = $wb.Worksheets.Item("Germany")
$wsGermany.Cells.Select()

# Create the new "Swiss" sheet by copying the "Czech" sheet (same layout/styles/merges)
# and placing it right after "Czech" (i.e. at the end of the tab strip).
$wsCzech = $wb.Worksheets.Item("Czech")
$wsCzech.Copy([System.Reflection.Missing]::Value, $wsCzech)
$wsSwiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSwiss.Name = "Swiss"

# Populate the Switzerland-specific values
$wsSwiss.Range("B2").Value = "Switzerland Market"
$wsSwiss.Range("B4").Value = "NGC-3476/T2653/T2654/T2656"

# Match the active-cell/selection recorded for the new sheet
$wsSwiss.Range("B5").Select()
